$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = -0
$ws.Range("B2").Value = -0.1207136776600321
$ws.Range("D2").Value = 0.2423024455985499
$ws.Range("E2").Value = 0.005374176121488893
$ws.Range("F2").Value = -0
$ws.Range("G2").Value = 0
$ws.Range("I2").Value = -0
$ws.Range("J2").Value = -0
$ws.Range("K2").Value = -0.07541815276151176
$ws.Range("L2").Value = -0
$ws.Range("M2").Value = 0.2149985903168879
$ws.Range("N2").Value = 0.02294959987369615
$ws.Range("R2").Value = -0
$ws.Range("S2").Value = 0
$ws.Range("T2").Value = -0.1063339094882755
$ws.Range("V2").Value = 0.01394680953777648
$ws.Range("W2").Value = -0.01137358646709177
$ws.Range("Z2").Value = -0
$ws.Range("AB2").Value = 0
$ws.Range("AC2").Value = 0.01658479900506774
$ws.Range("AD2").Value = 0
$ws.Range("AE2").Value = 0.01208351808461532
$ws.Range("AF2").Value = -0.0004440765254445339
$ws.Range("AG2").Value = -0
$ws.Range("AI2").Value = -0
$ws.Range("AJ2").Value = 0
$ws.Range("AK2").Value = -0
$ws.Range("AL2").Value = -0.03554061054620306
$ws.Range("AM2").Value = 0
$ws.Range("AN2").Value = 0.04021729487905552
$ws.Range("AO2").Value = 0.07759915751275892
$ws.Range("AQ2").Value = 0
$ws.Range("AR2").Value = -0
$ws.Range("AT2").Value = 0
$ws.Range("AU2").Value = -0.1767592602726039
$ws.Range("AW2").Value = 0.07839264346614182
$ws.Range("AX2").Value = -0.0108300864088781
$ws.Range("AY2").Value = -0
$ws.Range("BB2").Value = -0
$ws.Range("BC2").Value = -0
$ws.Range("BD2").Value = -0.01477925423043125
$ws.Range("BF2").Value = 0.09076894492908641
$ws.Range("BG2").Value = 0.02907409353348358
$ws.Range("BI2").Value = 0
$ws.Range("BJ2").Value = -0
$ws.Range("BL2").Value = 0
$ws.Range("BM2").Value = 0.0300776765823109
$ws.Range("BO2").Value = -0.05324375875133324
$ws.Range("BP2").Value = -0.0816346210201504
$ws.Range("BU2").Value = 0
$ws.Range("BV2").Value = -0.07852575589901419
$ws.Range("BW2").Value = 0
$ws.Range("BX2").Value = 0.01899427025985643
$ws.Range("BY2").Value = -0.02409089098383337
$ws.Range("BZ2").Value = -0
$ws.Range("CD2").Value = -0
$ws.Range("CE2").Value = 0.0351828300481634
$ws.Range("CF2").Value = -0
$ws.Range("CG2").Value = -0.03740581133404869
$ws.Range("CH2").Value = 0.01959344589137578
$ws.Range("CI2").Value = 0
$ws.Range("CJ2").Value = -0
$ws.Range("CL2").Value = 0
$ws.Range("CM2").Value = -0
$ws.Range("CN2").Value = -0.0109574491941789
$ws.Range("CO2").Value = -0
$ws.Range("CP2").Value = 0.02695512644018117
$ws.Range("CQ2").Value = 0.04735313135584056
$ws.Range("CT2").Value = 0
$ws.Range("CU2").Value = -0
$ws.Range("CV2").Value = -0
$ws.Range("CW2").Value = 0.04972248908286744
$ws.Range("CY2").Value = -0.04294037083382247
$ws.Range("CZ2").Value = 0.01300290137226959
$ws.Range("DE2").Value = -0
$ws.Range("DF2").Value = 0.03367236878889769
$ws.Range("DH2").Value = 0.01883731363827525
$ws.Range("DI2").Value = 0.04878514119211853
$ws.Range("DJ2").Value = 0
$ws.Range("DK2").Value = -0
$ws.Range("DL2").Value = -0
$ws.Range("DN2").Value = 0
$ws.Range("DO2").Value = -0.01724963747544638
$ws.Range("DQ2").Value = 0.04923827118775025
$ws.Range("DR2").Value = -0.01741754025609108
$ws.Range("DS2").Value = -0
$ws.Range("DT2").Value = 0
$ws.Range("DU2").Value = -0
$ws.Range("DV2").Value = -0
$ws.Range("DW2").Value = 0
$ws.Range("DX2").Value = -0.06232683367924737
$ws.Range("DY2").Value = -0
$ws.Range("DZ2").Value = -0.008166624026259794
$ws.Range("EA2").Value = -0.03432712046390961
$ws.Range("EB2").Value = 0
$ws.Range("ED2").Value = 0
$ws.Range("EF2").Value = -0
$ws.Range("EG2").Value = 0.04419106693749611
$ws.Range("EI2").Value = 0.0938691243097028
$ws.Range("EJ2").Value = -0.03337327474115563
$ws.Range("EO2").Value = 0
$ws.Range("EP2").Value = 0.05524633990337205
$ws.Range("EQ2").Value = 0
$ws.Range("ER2").Value = -0.05458425257331163
$ws.Range("ES2").Value = 0.0181611645046857
$ws.Range("ET2").Value = 0
$ws.Range("EV2").Value = 0
$ws.Range("EX2").Value = 0
$ws.Range("EY2").Value = 0.04501369381718781
$ws.Range("EZ2").Value = 0
$ws.Range("FA2").Value = -0.03439071615201791
$ws.Range("FB2").Value = 0.01959492007179191
$ws.Range("FD2").Value = -0
$ws.Range("FF2").Value = -0
$ws.Range("FG2").Value = -0
$ws.Range("FH2").Value = -0.01283727506370472
$ws.Range("FJ2").Value = -0.006750623608650416
$ws.Range("FK2").Value = 0.01209768093104599
$ws.Range("FL2").Value = -0
$ws.Range("FM2").Value = 0
$ws.Range("FP2").Value = -0
$ws.Range("FQ2").Value = -0.01474334462705824
$ws.Range("FR2").Value = -0
$ws.Range("FS2").Value = 0.007985635772036061
$ws.Range("FT2").Value = -0.00508330136040864
$ws.Range("FV2").Value = -0
$ws.Range("FY2").Value = 0
$ws.Range("FZ2").Value = -0.03753255017633956
$ws.Range("GB2").Value = 0.02009416768528361
$ws.Range("GD2").Value = 0
$ws.Range("GE2").Value = -0
